$d = $word.ActiveDocument

# Locate the paragraph that currently ends the body text ("Now testing
# after adding support for github hosts other than github.com (again, 2,
# 3)") so a brand-new paragraph can be inserted right after it (still
# ahead of the closing bookmark).
$target = $d.Content
$target.Find.ClearFormatting()
$found = $target.Find.Execute( `
    "Now testing after adding support for github hosts other than github.com (again, 2, 3)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the anchor paragraph to insert after."
}

# Collapse to a zero-length point right after the matched text and
# insert a fully-formed new paragraph (style BodyText) there, using raw
# WordprocessingML so the run/text node come out exactly like Word would
# author them (xml:space="preserve" included).
$insertionPoint = $d.Range($target.End, $target.End)
$insertionPoint.InsertXML( `
    "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
    "<w:pPr><w:pStyle w:val='BodyText'/></w:pPr>" + `
    "<w:r><w:t xml:space='preserve'>Testing that the check for a branch works… if it exist</w:t></w:r>" + `
    "</w:p>")
